$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# rows to append: row index, date serial, B, C, D
$data = @(
    @(230, 44304, 1, 4, 101.7293997965412),
    @(231, 44305, 0, 3, 76.2970498474059),
    @(232, 44306, 0, 2, 50.8646998982706),
    @(233, 44307, 0, 1, 25.4323499491353)
)

# template cell (A229) carries the date style used throughout column A
$styleSource = $ws.Cells.Item(229, 1)

foreach ($row in $data) {
    $r = $row[0]

    $cellA = $ws.Cells.Item($r, 1)
    $styleSource.Copy() | Out-Null
    $cellA.PasteSpecial(-4122) | Out-Null
    $cellA.Value = $row[1]

    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}

$excel.CutCopyMode = $false
